$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text content looks like a plain number (e.g. "1.00", "302.50")
# must be forced to Text format first, otherwise Excel auto-converts them to
# numeric values and the literal formatting (trailing zeros, etc.) is lost.
# The number format is restored to General afterwards so the cell formatting
# is unaffected by this workaround.

$ws.Range("D2").Value = '42.054.46'
$ws.Range("E2").Value = '  +5.37%  '
$ws.Range("D3").Value = '2.271.39'
$ws.Range("E3").Value = '  +2.43%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").NumberFormat = "General"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '302.50'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +3.45%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.11'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = '  +6.50%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.534'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  +4.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = '  -0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.488'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = '  +4.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '32.97'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +7.29%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.80'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  +9.69%  '
$ws.Range("E12").Value = '  +2.80%  '
$ws.Range("E13").Value = '  +3.14%  '
$ws.Range("E14").Value = '  +3.81%  '
$ws.Range("D15").Value = '2.616.82'
$ws.Range("E15").Value = '  +2.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.23'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  +3.54%  '
$ws.Range("D17").Value = '2.266.51'
$ws.Range("E17").Value = '  +2.17%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.758'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  +3.53%  '
$ws.Range("D19").Value = '41.948.19'
$ws.Range("E19").Value = '  +5.24%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.25'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  +9.11%  '
$ws.Range("E21").Value = '  +2.78%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.95'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  +3.36%  '
$ws.Range("E23").Value = '  +2.72%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '242.65'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  +2.30%  '
$ws.Range("E25").Value = '  +5.62%  '
$ws.Range("E26").Value = '  -0.07%  '
$ws.Range("E27").Value = '  +4.71%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '23.98'
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = '  +2.43%  '
$ws.Range("E29").Value = '  +1.95%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.72'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  +5.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.15'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  +6.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '158.06'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = '  +0.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = '  +0.03%  '
$ws.Range("E34").Value = '  +4.76%  '
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.09'
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = '  +5.55%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0743'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  +4.46%  '
$ws.Range("E37").Value = '  +3.11%  '
$ws.Range("E38").Value = '  +6.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '16.61'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = '  +8.40%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.80'
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = '  +5.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.99'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +6.71%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.17'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  +13.51%  '
$ws.Range("D44").Value = '2.050.28'
$ws.Range("E44").Value = '  -2.97%  '
$ws.Range("E45").Value = '  +4.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.10'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = '  +1.54%  '
$ws.Range("E47").Value = '  +7.58%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.00'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = '  -4.28%  '
$ws.Range("D49").Value = '2.491.45'
$ws.Range("E49").Value = '  +2.44%  '
$ws.Range("E50").Value = '  +2.17%  '
$ws.Range("E51").Value = '  +4.82%  '
